$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the updated cells keep their original text representation
# (these columns store formatted numeric/percentage strings as text,
# not real numbers, so force text number-format before assigning).
$ws.Range("D2:E51").NumberFormat = "@"

# --- Updated price (D) / volume change (E) values ---
$ws.Range("D2").Value = "29.840.58"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.888.25"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "0.7484"
$ws.Range("E5").Value = "  -4.39%  "
$ws.Range("D6").Value = "242.34"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("D8").Value = "0.3132"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "25.33"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").Value = "0.07133"
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("D11").Value = "0.08503"
$ws.Range("E11").Value = "  +5.08%  "
$ws.Range("D12").Value = "0.7609"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").Value = "1.881.89"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "5.366"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "93.41"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "6.171"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "29.692.88"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "13.71"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("D19").Value = "243.51"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("D20").Value = "0.000007804"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "2.151.80"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").Value = "0.9993"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "7.979"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "0.1591"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "9.359"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").Value = "18.75"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").Value = "2.028"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").Value = "1.499"
$ws.Range("E30").Value = "  +4.14%  "
$ws.Range("D31").Value = "1.538"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "4.507"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").Value = "4.120"
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D34").Value = "0.05416"
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("D35").Value = "1.240"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "0.7459"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "2.706"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").Value = "0.01943"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("D41").Value = "0.4457"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").Value = "1.089.38"
$ws.Range("E43").Value = "  -4.26%  "
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("D45").Value = "0.8535"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "1.0000"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D49").Value = "1.862"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").Value = "3.047"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("D51").Value = "2.043.13"
$ws.Range("E51").Value = "  +0.41%  "

# --- Rows 47 and 48 swapped: Aptos now ranks above Quant ---
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.699"
$ws.Range("E47").Value = "  +2.38%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "102.23"
$ws.Range("E48").Value = "  -0.28%  "
